# Bomporto.xlsx: fill in the previously-missing 2021-02-08 data point
# (date serial 44235) and extend the daily series through 2021-03-02
# (44257). Every row from the insertion point onward is renumbered and
# its 7-day rolling sum / per-100k-abitanti figures (cols C/D) recompute.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 90-112: write the post-insert values directly (equivalent to
# inserting a row at 93 and recalculating, but avoids COM row-insert
# style side effects since every target row already carries the date
# style from the original sheet).
$ws.Range("A90").Value = 44232
$ws.Range("B90").Value = 9
$ws.Range("C90").Value = 38
$ws.Range("D90").Value = 377.3584905660377

$ws.Range("A91").Value = 44233
$ws.Range("B91").Value = 9
$ws.Range("C91").Value = 33
$ws.Range("D91").Value = 327.7060575968222

$ws.Range("A92").Value = 44234
$ws.Range("B92").Value = 4
$ws.Range("C92").Value = 33
$ws.Range("D92").Value = 327.7060575968222

$ws.Range("A93").Value = 44235
$ws.Range("B93").Value = 4
$ws.Range("C93").Value = 30
$ws.Range("D93").Value = 297.914597815293

$ws.Range("A94").Value = 44236
$ws.Range("B94").Value = 3
$ws.Range("C94").Value = 28
$ws.Range("D94").Value = 278.0536246276067

$ws.Range("A95").Value = 44237
$ws.Range("B95").Value = 0
$ws.Range("C95").Value = 22
$ws.Range("D95").Value = 218.4707050645482

$ws.Range("A96").Value = 44238
$ws.Range("B96").Value = 1
$ws.Range("C96").Value = 19
$ws.Range("D96").Value = 188.6792452830189

$ws.Range("A97").Value = 44239
$ws.Range("B97").Value = 7
$ws.Range("C97").Value = 19
$ws.Range("D97").Value = 188.6792452830189

$ws.Range("A98").Value = 44240
$ws.Range("B98").Value = 3
$ws.Range("C98").Value = 17
$ws.Range("D98").Value = 168.8182720953327

$ws.Range("A99").Value = 44241
$ws.Range("B99").Value = 1
$ws.Range("C99").Value = 17
$ws.Range("D99").Value = 168.8182720953327

$ws.Range("A100").Value = 44242
$ws.Range("B100").Value = 4
$ws.Range("C100").Value = 17
$ws.Range("D100").Value = 168.8182720953327

$ws.Range("A101").Value = 44243
$ws.Range("B101").Value = 1
$ws.Range("C101").Value = 12
$ws.Range("D101").Value = 119.1658391261172

$ws.Range("A102").Value = 44244
$ws.Range("B102").Value = 0
$ws.Range("C102").Value = 10
$ws.Range("D102").Value = 99.30486593843098

$ws.Range("A103").Value = 44245
$ws.Range("B103").Value = 1
$ws.Range("C103").Value = 14
$ws.Range("D103").Value = 139.0268123138034

$ws.Range("A104").Value = 44246
$ws.Range("B104").Value = 2
$ws.Range("C104").Value = 14
$ws.Range("D104").Value = 139.0268123138034

$ws.Range("A105").Value = 44247
$ws.Range("B105").Value = 1
$ws.Range("C105").Value = 19
$ws.Range("D105").Value = 188.6792452830189

$ws.Range("A106").Value = 44248
$ws.Range("B106").Value = 5
$ws.Range("C106").Value = 20
$ws.Range("D106").Value = 198.609731876862

$ws.Range("A107").Value = 44249
$ws.Range("B107").Value = 4
$ws.Range("C107").Value = 27
$ws.Range("D107").Value = 268.1231380337637

$ws.Range("A108").Value = 44250
$ws.Range("B108").Value = 6
$ws.Range("C108").Value = 36
$ws.Range("D108").Value = 357.4975173783516

$ws.Range("A109").Value = 44251
$ws.Range("B109").Value = 1
$ws.Range("C109").Value = 39
$ws.Range("D109").Value = 387.2889771598809

$ws.Range("A110").Value = 44252
$ws.Range("B110").Value = 8
$ws.Range("C110").Value = 45
$ws.Range("D110").Value = 446.8718967229395

$ws.Range("A111").Value = 44253
$ws.Range("B111").Value = 11
$ws.Range("C111").Value = 48
$ws.Range("D111").Value = 476.6633565044688

$ws.Range("A112").Value = 44254
$ws.Range("B112").Value = 4
$ws.Range("C112").Value = 50
$ws.Range("D112").Value = 496.5243296921549

# Row 113 now represents 2021-02-28 (44255); its 7-day window is still
# incomplete so C/D stay blank, same as the other not-yet-complete rows.
$ws.Range("A113").Value = 44255
$ws.Range("B113").Value = 11

# Rows 114-115 are brand new trailing rows (2021-03-01 and 2021-03-02).
# Clone the A-column formatting from the row above (thin border, bold,
# centered, custom date numFmt) so the new date cells match the rest of
# the column; C/D are left blank like every other not-yet-complete row.
$ws.Range("A113").Copy()
$ws.Range("A114").PasteSpecial(-4122)
$ws.Range("A114").Value = 44256
$ws.Range("B114").Value = 7

$ws.Range("A114").Copy()
$ws.Range("A115").PasteSpecial(-4122)
$ws.Range("A115").Value = 44257
$ws.Range("B115").Value = 8

$excel.CutCopyMode = $false
